# Auto update stock data
# Updates Date_1 (column A) from 2025/12/22 to 2025/12/23 for each company row,
# and refreshes the EBITDA (column B) values where applicable.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
}

# Row 2 - Alcoa
Set-TextValue $ws.Range("A2") "2025/12/23"
Set-TextValue $ws.Range("B2") "6.68"

# Row 8 - Rio Tinto
Set-TextValue $ws.Range("A8") "2025/12/23"
Set-TextValue $ws.Range("B8") "8.43"

# Row 14 - Norsk Hydro
Set-TextValue $ws.Range("A14") "2025/12/23"
Set-TextValue $ws.Range("B14") "3.00"

# Row 20 - Reliance
Set-TextValue $ws.Range("A20") "2025/12/23"
Set-TextValue $ws.Range("B20") "12.95"

# Row 26 - Kaiser
Set-TextValue $ws.Range("A26") "2025/12/23"
Set-TextValue $ws.Range("B26") "11.21"

# Row 32 - Ryerson
Set-TextValue $ws.Range("A32") "2025/12/23"
Set-TextValue $ws.Range("B32") "27.79"

# Row 38 - Alro Steel (date only; EBITDA unchanged)
Set-TextValue $ws.Range("A38") "2025/12/23"

# Row 44 - Ultra
Set-TextValue $ws.Range("A44") "2025/12/23"
Set-TextValue $ws.Range("B44") "11.18"

# Row 50 - Benchmark
Set-TextValue $ws.Range("A50") "2025/12/23"
Set-TextValue $ws.Range("B50") "11.53"

# Row 56 - Celestica
Set-TextValue $ws.Range("A56") "2025/12/23"
Set-TextValue $ws.Range("B56") "32.15"

# Row 62 - Jabil
Set-TextValue $ws.Range("A62") "2025/12/23"
Set-TextValue $ws.Range("B62") "11.47"

# Row 68 - Flex
Set-TextValue $ws.Range("A68") "2025/12/23"
Set-TextValue $ws.Range("B68") "13.26"

# Row 74 - MKS
Set-TextValue $ws.Range("A74") "2025/12/23"
Set-TextValue $ws.Range("B74") "16.59"
